$wb = $excel.ActiveWorkbook

# Update the bank details sheet (financialBankDetails) with new bank info
$wsBank = $wb.Worksheets.Item("financialBankDetails")
$wsBank.Range("B2").Value = "STATE BANK OF INDIA-SBI Tresury Branch, Kurnool"
$wsBank.Range("C2").Value = "4502106--844810206002--STATE BANK OF INDIA"

# Make financialBankDetails the active sheet/tab
$wsBank.Activate()

$wb.Save()
